# Auto-generated data refresh for Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10: A Jawbreaking Weapon of Staggering Weight
$ws.Range("H10").Value = 7833
$ws.Range("J10").Value = 7833
$ws.Range("L10").Value = 7833
$ws.Range("N10").Value = -8419

# Row 98: The Dotted Line
$ws.Range("H98").Value = 7772.3213
$ws.Range("I98").Value = 7865.24
$ws.Range("J98").Value = 6998
$ws.Range("K98").Value = 7865.24
$ws.Range("L98").Value = 6998
$ws.Range("M98").Value = -6367.24
$ws.Range("N98").Value = -9994

# Row 122: Wishful Inking
$ws.Range("H122").Value = 7772.3213
$ws.Range("I122").Value = 7865.24
$ws.Range("J122").Value = 6998
$ws.Range("K122").Value = 23595.72
$ws.Range("L122").Value = 20994
$ws.Range("M122").Value = -21145.72
$ws.Range("N122").Value = -25894

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3779.323
$ws.Range("I32").Value = 3682.125
$ws.Range("K32").Value = 3682.125
$ws.Range("M32").Value = -3395.125

# Row 92: Mail It In
$ws.Range("H92").Value = 44996
$ws.Range("J92").Value = 44996
$ws.Range("L92").Value = 44996
$ws.Range("N92").Value = -49988

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal
$ws.Range("H94").Value = 4666.8335
$ws.Range("I94").Value = 4000.25
$ws.Range("K94").Value = 4000.25
$ws.Range("M94").Value = -3549.25

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 7703.25
$ws.Range("I134").Value = 7495.2
$ws.Range("K134").Value = 22485.6
$ws.Range("M134").Value = -19950.6

$ws = $wb.Worksheets.Item("CRP")
# Row 33: Tools for the Tools
$ws.Range("H33").Value = 5015.5
$ws.Range("I33").Value = 5015.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 5015.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -4636.5
$ws.Range("N33").Value = $null

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 6999
$ws.Range("I62").Value = 6999
$ws.Range("K62").Value = 6999
$ws.Range("M62").Value = -6375

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 6999
$ws.Range("I65").Value = 6999
$ws.Range("K65").Value = 34995
$ws.Range("M65").Value = -31875

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 5544.4443
$ws.Range("I122").Value = 5300
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 15900
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -13450
$ws.Range("N122").Value = -27400

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 3124.9473
$ws.Range("I132").Value = 2836.4614
$ws.Range("K132").Value = 8509.3842
$ws.Range("M132").Value = -5979.3842

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa
$ws.Range("H3").Value = 3241.125
$ws.Range("I3").Value = 3241.125
$ws.Range("K3").Value = 9723.375
$ws.Range("M3").Value = -9611.375

# Row 81: It Goes Down Smoothly
$ws.Range("H81").Value = 8795
$ws.Range("I81").Value = 453.66666
$ws.Range("K81").Value = 1360.99998
$ws.Range("M81").Value = -237.9999800000001

# Row 84: Quenching the Flame (L)
$ws.Range("H84").Value = 8795
$ws.Range("I84").Value = 453.66666
$ws.Range("K84").Value = 4082.99994
$ws.Range("M84").Value = 1533.00006

# Row 96: Hunger Is No Game
$ws.Range("H96").Value = 7622.25
$ws.Range("J96").Value = 9996.333000000001
$ws.Range("L96").Value = 29988.999
$ws.Range("N96").Value = -34106.999

# Row 99: A Shorlonging for the Familiar
$ws.Range("H99").Value = 1277.6666

# Row 128: A Historical Flavor
$ws.Range("H128").Value = 520000.12
$ws.Range("I128").Value = 520000.12
$ws.Range("K128").Value = 1560000.36
$ws.Range("M128").Value = -1555020.36

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 1495
$ws.Range("I80").Value = 1495
$ws.Range("K80").Value = 1495
$ws.Range("M80").Value = -497

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 1495
$ws.Range("I83").Value = 1495
$ws.Range("K83").Value = 7475
$ws.Range("M83").Value = -2483

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 3414.3
$ws.Range("I97").Value = 835.0714
$ws.Range("J97").Value = 9432.5
$ws.Range("K97").Value = 835.0714
$ws.Range("L97").Value = 9432.5
$ws.Range("M97").Value = -339.0714
$ws.Range("N97").Value = -10424.5

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 6732.6
$ws.Range("I102").Value = 5814
$ws.Range("K102").Value = 5814
$ws.Range("M102").Value = -4192

# Row 133: Pendulums of Our Own
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -95120

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1668.5
$ws.Range("I22").Value = 1084
$ws.Range("K22").Value = 1084
$ws.Range("M22").Value = -789

# Row 27: Fire and Hide
$ws.Range("H27").Value = 1668.5
$ws.Range("I27").Value = 1084
$ws.Range("K27").Value = 1084
$ws.Range("M27").Value = -977

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 682.5
$ws.Range("I55").Value = 846.44446
$ws.Range("J55").Value = 518.55554
$ws.Range("K55").Value = 846.44446
$ws.Range("L55").Value = 518.55554
$ws.Range("M55").Value = -673.44446
$ws.Range("N55").Value = -864.55554

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 9879.571
$ws.Range("I68").Value = 8313.223
$ws.Range("K68").Value = 8313.223
$ws.Range("M68").Value = -7564.223

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 9879.571
$ws.Range("I71").Value = 8313.223
$ws.Range("K71").Value = 41566.115
$ws.Range("M71").Value = -37822.115

# Row 80: Don't Sweat the Small Fry
$ws.Range("H80").Value = 58000
$ws.Range("J80").Value = 58000
$ws.Range("L80").Value = 58000
$ws.Range("N80").Value = -60246

# Row 81: I Need Your Glove Tonight
$ws.Range("H81").Value = 80000
$ws.Range("J81").Value = 80000
$ws.Range("L81").Value = 80000
$ws.Range("N81").Value = -81996

# Row 83: It's All in the Wrists (L)
$ws.Range("H83").Value = 58000
$ws.Range("J83").Value = 58000
$ws.Range("L83").Value = 174000
$ws.Range("N83").Value = -185232

# Row 84: Halonic Drake Handlers (L)
$ws.Range("H84").Value = 80000
$ws.Range("J84").Value = 80000
$ws.Range("L84").Value = 240000
$ws.Range("N84").Value = -249984

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 7292.273
$ws.Range("J132").Value = 5326.75
$ws.Range("L132").Value = 15980.25
$ws.Range("N132").Value = -21040.25

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table
$ws.Range("H113").Value = 404.86365
$ws.Range("I113").Value = 405.57144
$ws.Range("K113").Value = 1216.71432
$ws.Range("M113").Value = 953.28568

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1699.2222
$ws.Range("I126").Value = 1484.2858
$ws.Range("J126").Value = 2451.5
$ws.Range("K126").Value = 4452.857400000001
$ws.Range("L126").Value = 7354.5
$ws.Range("M126").Value = -1982.857400000001
$ws.Range("N126").Value = -12294.5

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 7340.4375
$ws.Range("I132").Value = 5957.4614
$ws.Range("K132").Value = 17872.3842
$ws.Range("M132").Value = -15342.3842
